$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking price cells to remain text so trailing zeros / exact digits are preserved
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates from the diff
$ws.Range("D2").Value = '38.963.14'
$ws.Range("E2").Value = '  -3.67%  '
$ws.Range("D3").Value = '2.212.64'
$ws.Range("E3").Value = '  -6.22%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '295.01'
$ws.Range("E5").Value = '  -4.62%  '
$ws.Range("D6").Value = '81.84'
$ws.Range("E6").Value = '  -4.78%  '
$ws.Range("D7").Value = '0.507'
$ws.Range("E7").Value = '  -3.16%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -5.34%  '
$ws.Range("E10").Value = '  -7.35%  '
$ws.Range("D11").Value = '28.53'
$ws.Range("E11").Value = '  -5.75%  '
$ws.Range("D12").Value = '46.45'
$ws.Range("E12").Value = '  -11.43%  '
$ws.Range("E13").Value = '  -1.86%  '
$ws.Range("D14").Value = '2.552.67'
$ws.Range("E14").Value = '  -6.65%  '
$ws.Range("D15").Value = '6.15'
$ws.Range("E15").Value = '  -5.35%  '
$ws.Range("D16").Value = '13.95'
$ws.Range("E16").Value = '  -5.97%  '
$ws.Range("D17").Value = '2.215.64'
$ws.Range("E17").Value = '  -6.40%  '
$ws.Range("E18").Value = '  -5.34%  '
$ws.Range("D19").Value = '38.879.07'
$ws.Range("E19").Value = '  -3.89%  '
$ws.Range("E20").Value = '  -4.10%  '
$ws.Range("E21").Value = '  -6.43%  '
$ws.Range("D22").Value = '64.49'
$ws.Range("E22").Value = '  -5.48%  '
$ws.Range("D23").Value = '10.03'
$ws.Range("E23").Value = '  -5.58%  '
$ws.Range("D24").Value = '225.35'
$ws.Range("E24").Value = '  -2.67%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -9.32%  '
$ws.Range("D27").Value = '1.75'
$ws.Range("E27").Value = '  -2.17%  '
$ws.Range("D28").Value = '22.37'
$ws.Range("E28").Value = '  -4.36%  '
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").Value = '9.00'
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("D31").Value = '147.36'
$ws.Range("E31").Value = '  -2.91%  '
$ws.Range("D32").Value = '31.52'
$ws.Range("E32").Value = '  -5.63%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("D34").Value = '4.80'
$ws.Range("E34").Value = '  -6.54%  '
$ws.Range("D35").Value = '0.0692'
$ws.Range("E35").Value = '  -4.04%  '
$ws.Range("E36").Value = '  -5.10%  '
$ws.Range("E37").Value = '  -3.42%  '
$ws.Range("D38").Value = '2.65'
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("D39").Value = '0.0944'
$ws.Range("E39").Value = '  -3.64%  '
$ws.Range("D40").Value = '14.62'
$ws.Range("E40").Value = '  -6.26%  '
$ws.Range("E41").Value = '  -4.50%  '
$ws.Range("E42").Value = '  -3.50%  '
$ws.Range("D43").Value = '1.899.44'
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("E44").Value = '  -3.96%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '8.97'
$ws.Range("E46").Value = '  -4.31%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '15.90'
$ws.Range("E47").Value = '  -8.81%  '
$ws.Range("D48").Value = '2.59'
$ws.Range("E48").Value = '  -3.31%  '
$ws.Range("D49").Value = '2.415.53'
$ws.Range("E49").Value = '  -6.99%  '
$ws.Range("D50").Value = '70.74'
$ws.Range("E50").Value = '  -1.62%  '
$ws.Range("D51").Value = '86.64'
$ws.Range("E51").Value = '  -6.15%  '
